$wb = $excel.ActiveWorkbook

# Update "Latest Handoff Datetime" (column D) for the 98acb8c9 dependency row (row 6)
# on both locale report sheets, reflecting a new handoff generated for that file.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-10 00:43:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-10 00:43:32"
